$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{}
$data[2] = @{ "E"=2; "G"=30.003349; "H"=60.006698; "I"=0.1788086146751209; "J"=0.1319456513375379; "K"=2; "M"=1.8337055; "N"=3.667411; "O"=0.007313502929307318; "P"=0.00521519558215984; "Q"=55.0173060797195; "R"=220.069224318878; "S"=0.00130771732721188; "T"=0.0006881223779407304 }
$data[3] = @{ "E"=2; "G"=30.003349; "H"=60.006698; "I"=0.1788086146751209; "J"=0.1319456513375379; "K"=3; "M"=85.28552633333334; "N"=255.856579; "O"=0.3401505567095424; "P"=0.3638376229081851; "Q"=2558.85141122769; "R"=15353.10846736614; "S"=0.06082184982620441; "T"=0.04800679213572199 }
$data[4] = @{ "E"=2; "G"=30.003349; "H"=60.006698; "I"=0.1788086146751209; "J"=0.1319456513375379; "K"=3; "M"=32.77748500000001; "N"=98.33245500000001; "O"=0.1307288616208147; "P"=0.1398324280804446; "Q"=983.4343217972652; "R"=5900.60593078359; "S"=0.02337544664447344; "T"=0.01845028080118368 }
$data[5] = @{ "E"=2; "G"=30.003349; "H"=60.006698; "I"=0.1788086146751209; "J"=0.1319456513375379; "K"=3; "M"=80.65247599999999; "N"=241.957428; "O"=0.3216722202566814; "P"=0.3440725104375696; "Q"=2419.844385142124; "R"=14519.06631085274; "S"=0.05751776408356755; "T"=0.04539887149702693 }
$data[6] = @{ "E"=2; "G"=30.003349; "H"=60.006698; "I"=0.1788086146751209; "J"=0.1319456513375379; "K"=3; "M"=3.043386333333333; "N"=9.130158999999999; "O"=0.01213816224243598; "P"=0.01298342751363752; "Q"=91.31178230083032; "R"=547.8706938049819; "S"=0.002170407975271836; "T"=0.001713106799880613 }
$data[7] = @{ "E"=2; "G"=30.003349; "H"=60.006698; "I"=0.1788086146751209; "J"=0.1319456513375379; "K"=2; "M"=47.136178; "N"=94.272356; "O"=0.1879966962412182; "P"=0.1340588154780033; "Q"=1414.243199060122; "R"=5656.972796240488; "S"=0.03361542881839173; "T"=0.01768847772578396 }
$data[8] = @{ "E"=3; "G"=65.63887766666667; "H"=196.916633; "I"=0.3911828904302115; "J"=0.4329898872352535; "K"=2; "M"=1.8337055; "N"=3.667411; "O"=0.007313502929307318; "P"=0.00521519558215984; "Q"=120.3623709911938; "R"=722.1742259471631; "S"=0.002860917215056256; "T"=0.002258126947029181 }
$data[9] = @{ "E"=3; "G"=65.63887766666667; "H"=196.916633; "I"=0.3911828904302115; "J"=0.4329898872352535; "K"=3; "M"=85.28552633333334; "N"=255.856579; "O"=0.3401505567095424; "P"=0.3638376229081851; "Q"=5598.046229730946; "R"=50382.41606757852; "S"=0.1330610779550844; "T"=0.1575380113149578 }
$data[10] = @{ "E"=3; "G"=65.63887766666667; "H"=196.916633; "I"=0.3911828904302115; "J"=0.4329898872352535; "K"=3; "M"=32.77748500000001; "N"=98.33245500000001; "O"=0.1307288616208147; "P"=0.1398324280804446; "Q"=2151.477328136002; "R"=19363.29595322402; "S"=0.05113889395148143; "T"=0.06054602726638339 }
$data[11] = @{ "E"=3; "G"=65.63887766666667; "H"=196.916633; "I"=0.3911828904302115; "J"=0.4329898872352535; "K"=3; "M"=80.65247599999999; "N"=241.957428; "O"=0.3216722202566814; "P"=0.3440725104375696; "Q"=5293.938005677769; "R"=47645.44205109993; "S"=0.1258326688911123; "T"=0.1489799174951138 }
$data[12] = @{ "E"=3; "G"=65.63887766666667; "H"=196.916633; "I"=0.3911828904302115; "J"=0.4329898872352535; "K"=3; "M"=3.043386333333333; "N"=9.130158999999999; "O"=0.01213816224243598; "P"=0.01298342751363752; "Q"=199.7644632260719; "R"=1797.880169034647; "S"=0.004748241390506965; "T"=0.005621692815056999 }
$data[13] = @{ "E"=3; "G"=65.63887766666667; "H"=196.916633; "I"=0.3911828904302115; "J"=0.4329898872352535; "K"=2; "M"=47.136178; "N"=94.272356; "O"=0.1879966962412182; "P"=0.1340588154780033; "Q"=3093.965821416225; "R"=18563.79492849735; "S"=0.07354109102697023; "T"=0.05804611139671232 }
$data[14] = @{ "E"=3; "G"=13.43183266666667; "H"=40.295498; "I"=0.0800486436255733; "J"=0.08860370436614352; "K"=2; "M"=1.8337055; "N"=3.667411; "O"=0.007313502929307318; "P"=0.00521519558215984; "Q"=24.63002543594633; "R"=147.780152615678; "S"=0.000585435989642708; "T"=0.0004620856475733082 }
$data[15] = @{ "E"=3; "G"=13.43183266666667; "H"=40.295498; "I"=0.0800486436255733; "J"=0.08860370436614352; "K"=3; "M"=85.28552633333334; "N"=255.856579; "O"=0.3401505567095424; "P"=0.3638376229081851; "Q"=1145.540918597927; "R"=10309.86826738134; "S"=0.02722859069308252; "T"=0.03223736117743724 }
$data[16] = @{ "E"=3; "G"=13.43183266666667; "H"=40.295498; "I"=0.0800486436255733; "J"=0.08860370436614352; "K"=3; "M"=32.77748500000001; "N"=98.33245500000001; "O"=0.1307288616208147; "P"=0.1398324280804446; "Q"=440.2616937541767; "R"=3962.355243787591; "S"=0.01046466805546148; "T"=0.01238967111843974 }
$data[17] = @{ "E"=3; "G"=13.43183266666667; "H"=40.295498; "I"=0.0800486436255733; "J"=0.08860370436614352; "K"=3; "M"=80.65247599999999; "N"=241.957428; "O"=0.3216722202566814; "P"=0.3440725104375696; "Q"=1083.310561784349; "R"=9749.795056059143; "S"=0.02574942492357401; "T"=0.03048609899532725 }
$data[18] = @{ "E"=3; "G"=13.43183266666667; "H"=40.295498; "I"=0.0800486436255733; "J"=0.08860370436614352; "K"=3; "M"=3.043386333333333; "N"=9.130158999999999; "O"=0.01213816224243598; "P"=0.01298342751363752; "Q"=40.87825596935355; "R"=367.904303724182; "S"=0.0009716434236141475; "T"=0.001150379773077593 }
$data[19] = @{ "E"=3; "G"=13.43183266666667; "H"=40.295498; "I"=0.0800486436255733; "J"=0.08860370436614352; "K"=2; "M"=47.136178; "N"=94.272356; "O"=0.1879966962412182; "P"=0.1340588154780033; "Q"=633.1252554422147; "R"=3798.751532653288; "S"=0.01504888054019844; "T"=0.01187810765428839 }
$data[20] = @{ "E"=3; "G"=16.42760466666667; "H"=49.282814; "I"=0.09790231193448497; "J"=0.108365452686244; "K"=2; "M"=1.8337055; "N"=3.667411; "O"=0.007313502929307318; "P"=0.00521519558215984; "Q"=30.12338902909233; "R"=180.740334174554; "S"=0.0007160088451188146; "T"=0.0005651470301080508 }
$data[21] = @{ "E"=3; "G"=16.42760466666667; "H"=49.282814; "I"=0.09790231193448497; "J"=0.108365452686244; "K"=3; "M"=85.28552633333334; "N"=255.856579; "O"=0.3401505567095424; "P"=0.3638376229081851; "Q"=1401.03691039259; "R"=12609.33219353331; "S"=0.03330152590766634; "T"=0.03942742871073242 }
$data[22] = @{ "E"=3; "G"=16.42760466666667; "H"=49.282814; "I"=0.09790231193448497; "J"=0.108365452686244; "K"=3; "M"=32.77748500000001; "N"=98.33245500000001; "O"=0.1307288616208147; "P"=0.1398324280804446; "Q"=538.4555655475967; "R"=4846.100089928371; "S"=0.01279865778924112; "T"=0.01515300436915403 }
$data[23] = @{ "E"=3; "G"=16.42760466666667; "H"=49.282814; "I"=0.09790231193448497; "J"=0.108365452686244; "K"=3; "M"=80.65247599999999; "N"=241.957428; "O"=0.3216722202566814; "P"=0.3440725104375696; "Q"=1324.926991115821; "R"=11924.34292004239; "S"=0.03149245404822798; "T"=0.03728557335045964 }
$data[24] = @{ "E"=3; "G"=16.42760466666667; "H"=49.282814; "I"=0.09790231193448497; "J"=0.108365452686244; "K"=3; "M"=3.043386333333333; "N"=9.130158999999999; "O"=0.01213816224243598; "P"=0.01298342751363752; "Q"=49.99554753193622; "R"=449.959927787426; "S"=0.001188354146170355; "T"=0.001406954999934365 }
$data[25] = @{ "E"=3; "G"=16.42760466666667; "H"=49.282814; "I"=0.09790231193448497; "J"=0.108365452686244; "K"=2; "M"=47.136178; "N"=94.272356; "O"=0.1879966962412182; "P"=0.1340588154780033; "Q"=774.3344976816308; "R"=4646.006986089785; "S"=0.01840531119806036; "T"=0.01452734422585548 }
$data[26] = @{ "E"=3; "G"=23.69336933333333; "H"=71.080108; "I"=0.1412035218961498; "J"=0.1562944047879878; "K"=2; "M"=1.8337055; "N"=3.667411; "O"=0.007313502929307318; "P"=0.00521519558215984; "Q"=43.44666166006466; "R"=260.679969960388; "S"=0.001032692371016002; "T"=0.0008151058893666158 }
$data[27] = @{ "E"=3; "G"=23.69336933333333; "H"=71.080108; "I"=0.1412035218961498; "J"=0.1562944047879878; "K"=3; "M"=85.28552633333334; "N"=255.856579; "O"=0.3401505567095424; "P"=0.3638376229081851; "Q"=2020.701474203393; "R"=18186.31326783053; "S"=0.04803045658232343; "T"=0.05686578471191115 }
$data[28] = @{ "E"=3; "G"=23.69336933333333; "H"=71.080108; "I"=0.1412035218961498; "J"=0.1562944047879878; "K"=3; "M"=32.77748500000001; "N"=98.33245500000001; "O"=0.1307288616208147; "P"=0.1398324280804446; "Q"=776.6090579227935; "R"=6989.481521305141; "S"=0.01845937567433344; "T"=0.0218550261168922 }
$data[29] = @{ "E"=3; "G"=23.69336933333333; "H"=71.080108; "I"=0.1412035218961498; "J"=0.1562944047879878; "K"=3; "M"=80.65247599999999; "N"=241.957428; "O"=0.3216722202566814; "P"=0.3440725104375696; "Q"=1910.928901515802; "R"=17198.36011364222; "S"=0.04542125039639745; "T"=0.05377660822274866 }
$data[30] = @{ "E"=3; "G"=23.69336933333333; "H"=71.080108; "I"=0.1412035218961498; "J"=0.1562944047879878; "K"=3; "M"=3.043386333333333; "N"=9.130158999999999; "O"=0.01213816224243598; "P"=0.01298342751363752; "Q"=72.10807641968576; "R"=648.9726877771719; "S"=0.001713951257978828; "T"=0.002029237075351961 }
$data[31] = @{ "E"=3; "G"=23.69336933333333; "H"=71.080108; "I"=0.1412035218961498; "J"=0.1562944047879878; "K"=2; "M"=47.136178; "N"=94.272356; "O"=0.1879966962412182; "P"=0.1340588154780033; "Q"=1116.814874315741; "R"=6700.889245894447; "S"=0.02654579561410069; "T"=0.02095264277171722 }
$data[32] = @{ "E"=2; "G"=18.6008475; "H"=37.201695; "I"=0.1108540174384595; "J"=0.08180089958683325; "K"=2; "M"=1.8337055; "N"=3.667411; "O"=0.007313502929307318; "P"=0.00521519558215984; "Q"=34.10847636541125; "R"=136.433905461645; "S"=0.0008107311812616578; "T"=0.0004266076901419534 }
$data[33] = @{ "E"=2; "G"=18.6008475; "H"=37.201695; "I"=0.1108540174384595; "J"=0.08180089958683325; "K"=3; "M"=85.28552633333334; "N"=255.856579; "O"=0.3401505567095424; "P"=0.3638376229081851; "Q"=1586.383069283568; "R"=9518.298415701405; "S"=0.03770705574518131; "T"=0.02976224485742455 }
$data[34] = @{ "E"=2; "G"=18.6008475; "H"=37.201695; "I"=0.1108540174384595; "J"=0.08180089958683325; "K"=3; "M"=32.77748500000001; "N"=98.33245500000001; "O"=0.1307288616208147; "P"=0.1398324280804446; "Q"=609.6889999185377; "R"=3658.133999511225; "S"=0.01449181950582374; "T"=0.01143841840839153 }
$data[35] = @{ "E"=2; "G"=18.6008475; "H"=37.201695; "I"=0.1108540174384595; "J"=0.08180089958683325; "K"=3; "M"=80.65247599999999; "N"=241.957428; "O"=0.3216722202566814; "P"=0.3440725104375696; "Q"=1500.20440657341; "R"=9001.226439440459; "S"=0.03565865791380213; "T"=0.02814544087689326 }
$data[36] = @{ "E"=2; "G"=18.6008475; "H"=37.201695; "I"=0.1108540174384595; "J"=0.08180089958683325; "K"=3; "M"=3.043386333333333; "N"=9.130158999999999; "O"=0.01213816224243598; "P"=0.01298342751363752; "Q"=56.60956506991749; "R"=339.657390419505; "S"=0.001345564048893848; "T"=0.001062056050335991 }
$data[37] = @{ "E"=2; "G"=18.6008475; "H"=37.201695; "I"=0.1108540174384595; "J"=0.08180089958683325; "K"=2; "M"=47.136178; "N"=94.272356; "O"=0.1879966962412182; "P"=0.1340588154780033; "Q"=876.7728587108551; "R"=3507.09143484342; "S"=0.02084018904349677; "T"=0.01096613170364596 }

foreach ($row in $data.Keys) {
    foreach ($col in $data[$row].Keys) {
        $ws.Range("$col$row").Value = $data[$row][$col]
    }
}